$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values for row 2 (angInc scenario)
$ws.Range("B2").Value = 87.443
$ws.Range("C2").Value = 0.602
$ws.Range("D2").Value = 0.0559

# Update input values for row 5 (lat scenario)
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 0.602
$ws.Range("D5").Value = 0.0559

# K10 holds a static copy of the K5 result (pasted as value)
$ws.Range("K10").Value = 89.017641452899696

# Update selection to match the diff (active cell K10)
$ws.Range("K10").Select()
